$wb = $excel.ActiveWorkbook

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6570.5
$ws.Range("I76").Value = 5662
$ws.Range("J76").Value = 7277.1113
$ws.Range("K76").Value = 5662
$ws.Range("L76").Value = 7277.1113
$ws.Range("M76").Value = -5347
$ws.Range("N76").Value = -7907.1113

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 6570.5
$ws.Range("I79").Value = 5662
$ws.Range("J79").Value = 7277.1113
$ws.Range("K79").Value = 5662
$ws.Range("L79").Value = 7277.1113
$ws.Range("M79").Value = -4570
$ws.Range("N79").Value = -9461.1113

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 896.7619
$ws.Range("I98").Value = 801.6842
$ws.Range("K98").Value = 801.6842
$ws.Range("M98").Value = 696.3158

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 896.7619
$ws.Range("I122").Value = 801.6842
$ws.Range("K122").Value = 2405.0526
$ws.Range("M122").Value = 44.94740000000002

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 795
$ws.Range("I125").Value = 907.3333
$ws.Range("K125").Value = 8165.9997
$ws.Range("M125").Value = -5705.9997

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2506.7856
$ws.Range("I137").Value = 1214.5714
$ws.Range("K137").Value = 3643.7142
$ws.Range("M137").Value = -1093.7142

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4405.651
$ws.Range("I138").Value = 1977.375
$ws.Range("J138").Value = 4758.8545
$ws.Range("K138").Value = 5932.125
$ws.Range("L138").Value = 14276.5635
$ws.Range("M138").Value = -792.125
$ws.Range("N138").Value = -24556.5635

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4032.2856
$ws.Range("I45").Value = 1621.3334
$ws.Range("K45").Value = 1621.3334
$ws.Range("M45").Value = -1244.3334

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4760.3335
$ws.Range("I74").Value = 1551.75
$ws.Range("J74").Value = 7327.2
$ws.Range("K74").Value = 1551.75
$ws.Range("L74").Value = 7327.2
$ws.Range("M74").Value = -677.75
$ws.Range("N74").Value = -9075.200000000001

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4760.3335
$ws.Range("I77").Value = 1551.75
$ws.Range("J77").Value = 7327.2
$ws.Range("K77").Value = 7758.75
$ws.Range("L77").Value = 36636
$ws.Range("M77").Value = -3390.75
$ws.Range("N77").Value = -45372

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1115522.4
$ws.Range("I122").Value = 2503675.2
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7511025.600000001
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -7508575.600000001
$ws.Range("N122").Value = -19900

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4788.8184
$ws.Range("I134").Value = 4788.8184
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14366.4552
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -11831.4552
$ws.Range("N134").ClearContents()

# CRP row 93
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 2392.3333
$ws.Range("I93").Value = 2392.3333
$ws.Range("K93").Value = 2392.3333
$ws.Range("M93").Value = -520.3332999999998

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1999
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1999
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1999
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2901

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13284.542
$ws.Range("I99").Value = 10305
$ws.Range("J99").Value = 16264.083
$ws.Range("K99").Value = 10305
$ws.Range("L99").Value = 16264.083
$ws.Range("M99").Value = -8807
$ws.Range("N99").Value = -19260.083

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 306.7143
$ws.Range("I107").Value = 234.15384
$ws.Range("K107").Value = 234.15384
$ws.Range("M107").Value = 1685.84616

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 13284.542
$ws.Range("I126").Value = 10305
$ws.Range("J126").Value = 16264.083
$ws.Range("K126").Value = 30915
$ws.Range("L126").Value = 48792.249
$ws.Range("M126").Value = -28445
$ws.Range("N126").Value = -53732.249

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5181.8
$ws.Range("I134").Value = 3672.5
$ws.Range("J134").Value = 6188
$ws.Range("K134").Value = 11017.5
$ws.Range("L134").Value = 18564
$ws.Range("M134").Value = -8482.5
$ws.Range("N134").Value = -23634

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2374.1667
$ws.Range("I39").Value = 1749
$ws.Range("K39").Value = 5247
$ws.Range("M39").Value = -4953

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 189450
$ws.Range("I55").Value = 250150
$ws.Range("J55").Value = 128750
$ws.Range("K55").Value = 750450
$ws.Range("L55").Value = 386250
$ws.Range("M55").Value = -750273
$ws.Range("N55").Value = -386604

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7998
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7998
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 849930.0600000001
$ws.Range("I122").Value = 128263.375
$ws.Range("K122").Value = 384790.125
$ws.Range("M122").Value = -382340.125

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2250.75
$ws.Range("I40").Value = 2250.75
$ws.Range("K40").Value = 2250.75
$ws.Range("M40").Value = -2114.75

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2168.6086
$ws.Range("I61").Value = 1843.8
$ws.Range("K61").Value = 1843.8
$ws.Range("M61").Value = -1641.8

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2799.2
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2799.2
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2799.2
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4297.2

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2799.2
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2799.2
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 13996
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -21484

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2168.6086
$ws.Range("I113").Value = 1843.8
$ws.Range("K113").Value = 1843.8
$ws.Range("M113").Value = 326.2

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2998.5
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 2997
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 8991
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -13891

# WVR row 60
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 99997.664
$ws.Range("J60").Value = 99999.5
$ws.Range("L60").Value = 99999.5
$ws.Range("N60").Value = -101643.5

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6642.143
$ws.Range("I62").Value = 1997.5
$ws.Range("J62").Value = 8500
$ws.Range("K62").Value = 1997.5
$ws.Range("L62").Value = 8500
$ws.Range("M62").Value = -1373.5
$ws.Range("N62").Value = -9748

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 6642.143
$ws.Range("I65").Value = 1997.5
$ws.Range("J65").Value = 8500
$ws.Range("K65").Value = 9987.5
$ws.Range("L65").Value = 42500
$ws.Range("M65").Value = -6867.5
$ws.Range("N65").Value = -48740

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1712.625
$ws.Range("I100").Value = 1528.7142
$ws.Range("K100").Value = 3057.4284
$ws.Range("M100").Value = -2516.4284

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 79590
$ws.Range("I126").Value = 92515.45
$ws.Range("J126").Value = 8500
$ws.Range("K126").Value = 277546.35
$ws.Range("L126").Value = 25500
$ws.Range("M126").Value = -275076.35
$ws.Range("N126").Value = -30440
